$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Scraper Success XML")

$newCreditCardXml = "<scrape-session>`r`n  <baseURL>www.xbox.com</baseURL>`r`n  <date>12/12/2014</date>`r`n  <time>13:50:00</time>`r`n  <datapair id=`"001`">`r`n    <text>Account Number</text>`r`n    <value>123456789</value>`r`n  </datapair>`r`n  <datapair id=`"002`">`r`n    <text>Account holder name</text>`r`n    <value>Jack Parcell</value>`r`n  </datapair>`r`n  <datapair id=`"003`">`r`n    <text>Statement date</text>`r`n    <value>12/12/2014</value>`r`n  </datapair>`r`n  <datapair id=`"004`">`r`n    <text>Statement number</text>`r`n    <value>1122</value>`r`n  </datapair>`r`n  <datapair id=`"005`">`r`n    <text>Statement month</text>`r`n    <value>2</value>`r`n  </datapair>`r`n  <datapair id=`"006`">`r`n    <text>Total due</text>`r`n    <value>R340</value>`r`n  </datapair>`r`n  <datapair id=`"007`">`r`n    <text>Due date</text>`r`n    <value>01/01/2015</value>`r`n  </datapair>`r`n  <datapair id=`"008`">`r`n    <text>Opening balance</text>`r`n    <value>R120</value>`r`n  </datapair>`r`n  <datapair id=`"009`">`r`n    <text>Closing balance</text>`r`n    <value>R123</value>`r`n  </datapair>`r`n  <datapair id=`"010`">`r`n    <text>Payment received</text>`r`n    <value>R40</value>`r`n  </datapair>`r`n  <datapair id=`"011`">`r`n    <text>New charges</text>`r`n    <value>R45</value>`r`n  </datapair>`r`n  <datapair id=`"012`">`r`n    <text>Deductions</text>`r`n    <value>R123</value>`r`n  </datapair>`r`n  <datapair id=`"013`">`r`n    <text>Discount</text>`r`n    <value>R456</value>`r`n  </datapair>`r`n  <datapair id=`"014`">`r`n    <text>VAT Amount</text>`r`n    <value>R123</value>`r`n  </datapair>`r`n  <datapair id=`"015`">`r`n    <text>Card type</text>`r`n    <value>Visa</value>`r`n  </datapair>`r`n  <datapair id=`"016`">`r`n    <text>Interest rate</text>`r`n    <value>12%</value>`r`n  </datapair>`r`n  <datapair id=`"017`">`r`n    <text>Credit limit</text>`r`n    <value>R20000</value>`r`n  </datapair>`r`n  <datapair id=`"018`">`r`n    <text>Credit available</text>`r`n    <value>R4500</value>`r`n  </datapair>`r`n  <datapair id=`"019`">`r`n    <text>Minimum amount due</text>`r`n    <value>R90</value>`r`n  </datapair>`r`n</scrape-session>"

$cell = $ws.Range("B2")
$cell.Value = $newCreditCardXml
$cell.NumberFormat = "General"
$cell.WrapText = $true
$ws.Rows.Item(2).RowHeight = 161.25

$ws.Activate() | Out-Null
$ws.Range("B2").Select() | Out-Null
